$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: reorder "Recorded By" list
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 3: reorder "Recorded By" list
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# Row 4: reorder + add recorder, update attendance count
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H4").Value = "86/251"

# Row 5: add recorder, update attendance count
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("H5").Value = "81/251"

# Row 7: reorder "Recorded By" list
$ws.Range("G7").Value = "Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"

# Row 10 summary (L10): average attendance % updated
$ws.Range("L10").Value = "24.9%"

# Row 15 summary (S15): average attendance % updated
$ws.Range("S15").Value = "24.9%"

# Row 28: reorder "Recorded By" list
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
